# Auto-update draw results: append the 2025-11-01 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 46

# Columns A (date-like) and C (all-digit) must stay TEXT, not be reinterpreted
# as a date serial / number, so they're entered with a leading apostrophe
# (Excel's standard "force text" quote-prefix), exactly like the rest of the
# column. B, D and E are not numeric-looking, so a plain Value assign already
# stores them as text.
$ws.Cells.Item($row, 1).Value = "'2025-11-01"
$ws.Cells.Item($row, 2).Value = "Pick 3"
$ws.Cells.Item($row, 3).Value = "'251101"
$ws.Cells.Item($row, 4).Value = "4-1-6"
$ws.Cells.Item($row, 5).Value = "2025-11-01T21:35:27.972+04:00"
